$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'college; university'
$ws.Cells.Item(2, 2).Value = 'だいがく'
$ws.Cells.Item(3, 1).Value = 'high school'
$ws.Cells.Item(3, 2).Value = 'こうこう'
$ws.Cells.Item(4, 1).Value = 'student'
$ws.Cells.Item(4, 2).Value = 'がくせい'
$ws.Cells.Item(5, 1).Value = 'college student'
$ws.Cells.Item(5, 2).Value = 'だいがくせい'
$ws.Cells.Item(6, 1).Value = 'international student'
$ws.Cells.Item(6, 2).Value = 'りゅうがくせい'
$ws.Cells.Item(7, 1).Value = 'teacher; professor'
$ws.Cells.Item(7, 2).Value = 'せんせい'
$ws.Cells.Item(8, 1).Value = '...year student'
$ws.Cells.Item(8, 2).Value = '～ねんせい'
$ws.Cells.Item(9, 1).Value = 'first-year student'
$ws.Cells.Item(9, 2).Value = 'いちねんせい'
$ws.Cells.Item(10, 1).Value = 'major'
$ws.Cells.Item(10, 2).Value = 'せんこう'
$ws.Cells.Item(11, 1).Value = '01:00'
$ws.Cells.Item(11, 2).Value = 'いちじ'
$ws.Cells.Item(12, 1).Value = '02:00'
$ws.Cells.Item(12, 2).Value = 'にじ'
$ws.Cells.Item(13, 1).Value = '03:00'
$ws.Cells.Item(13, 2).Value = 'さんじ'
$ws.Cells.Item(14, 1).Value = '04:00'
$ws.Cells.Item(14, 2).Value = 'よじ'
$ws.Cells.Item(15, 1).Value = '05:00'
$ws.Cells.Item(15, 2).Value = 'ごじ'
$ws.Cells.Item(16, 1).Value = '06:00'
$ws.Cells.Item(16, 2).Value = 'ろくじ'
$ws.Cells.Item(17, 1).Value = '07:00'
$ws.Cells.Item(17, 2).Value = 'しちじ'
$ws.Cells.Item(18, 1).Value = '08:00'
$ws.Cells.Item(18, 2).Value = 'はちじ'
$ws.Cells.Item(19, 1).Value = '09:00'
$ws.Cells.Item(19, 2).Value = 'くじ'
$ws.Cells.Item(20, 1).Value = '10:00'
$ws.Cells.Item(20, 2).Value = 'じゅうじ'
$ws.Cells.Item(21, 1).Value = '11:00'
$ws.Cells.Item(21, 2).Value = 'じゅういちじ'
$ws.Cells.Item(22, 1).Value = '12:00'
$ws.Cells.Item(22, 2).Value = 'じゅうにじ'
$ws.Cells.Item(23, 1).Value = '01:30'
$ws.Cells.Item(23, 2).Value = 'いちじはん'
$ws.Cells.Item(24, 1).Value = '00:01'
$ws.Cells.Item(24, 2).Value = 'いっぷん'
$ws.Cells.Item(25, 1).Value = '00:02'
$ws.Cells.Item(25, 2).Value = 'にふん'
$ws.Cells.Item(26, 1).Value = '00:03'
$ws.Cells.Item(26, 2).Value = 'さんぷん'
$ws.Cells.Item(27, 1).Value = '00:04'
$ws.Cells.Item(27, 2).Value = 'よんぷん'
$ws.Cells.Item(28, 1).Value = '00:05'
$ws.Cells.Item(28, 2).Value = 'ごふん'
$ws.Cells.Item(29, 1).Value = '00:06'
$ws.Cells.Item(29, 2).Value = 'ろっぷん'
$ws.Cells.Item(30, 1).Value = '00:07'
$ws.Cells.Item(30, 2).Value = 'ななふん'
$ws.Cells.Item(31, 1).Value = '00:08'
$ws.Cells.Item(31, 2).Value = 'はっぷん／はちふん'
$ws.Cells.Item(32, 1).Value = '00:09'
$ws.Cells.Item(32, 2).Value = 'きゅうふん'
$ws.Cells.Item(33, 1).Value = '00:10'
$ws.Cells.Item(33, 2).Value = 'じ(ゅ)っぷん'
$ws.Cells.Item(34, 1).Value = '00:11'
$ws.Cells.Item(34, 2).Value = 'じゅういっぷん'
$ws.Cells.Item(35, 1).Value = '00:12'
$ws.Cells.Item(35, 2).Value = 'じゅうにふん'
$ws.Cells.Item(36, 1).Value = '00:13'
$ws.Cells.Item(36, 2).Value = 'じゅうさんぷん'
$ws.Cells.Item(37, 1).Value = '00:14'
$ws.Cells.Item(37, 2).Value = 'じゅうよんぷん'
$ws.Cells.Item(38, 1).Value = '00:15'
$ws.Cells.Item(38, 2).Value = 'じゅうごふん'
$ws.Cells.Item(39, 1).Value = '00:16'
$ws.Cells.Item(39, 2).Value = 'じゅうろっぷん'
$ws.Cells.Item(40, 1).Value = '00:17'
$ws.Cells.Item(40, 2).Value = 'じゅうななふん'
$ws.Cells.Item(41, 1).Value = '00:18'
$ws.Cells.Item(41, 2).Value = 'じゅうはっぷん／じゅうはちふん'
$ws.Cells.Item(42, 1).Value = '00:19'
$ws.Cells.Item(42, 2).Value = 'じゅうきゅうふん'
$ws.Cells.Item(43, 1).Value = '00:20'
$ws.Cells.Item(43, 2).Value = 'にじ(ゅ)っぷん'
$ws.Cells.Item(44, 1).Value = '00:30'
$ws.Cells.Item(44, 2).Value = 'さんじ(ゅ)っぷん'
$ws.Cells.Item(45, 1).Value = 'I'
$ws.Cells.Item(45, 2).Value = 'わたし'
$ws.Cells.Item(46, 1).Value = 'friend'
$ws.Cells.Item(46, 2).Value = 'ともだち'
$ws.Cells.Item(47, 1).Value = 'Mr./Ms. ...'
$ws.Cells.Item(47, 2).Value = '～さん'
$ws.Cells.Item(48, 1).Value = '...people'
$ws.Cells.Item(48, 2).Value = '～じん'
$ws.Cells.Item(49, 1).Value = 'Japanese people'
$ws.Cells.Item(49, 2).Value = 'にほんじん'
$ws.Cells.Item(50, 1).Value = 'now'
$ws.Cells.Item(50, 2).Value = 'いま'
$ws.Cells.Item(51, 1).Value = 'A.M.'
$ws.Cells.Item(51, 2).Value = 'ごぜん'
$ws.Cells.Item(52, 1).Value = 'P.M.'
$ws.Cells.Item(52, 2).Value = 'ごご'
$ws.Cells.Item(53, 1).Value = 'o''clock'
$ws.Cells.Item(53, 2).Value = '～じ'
$ws.Cells.Item(54, 1).Value = 'one o''clock'
$ws.Cells.Item(54, 2).Value = 'いちじ'
$ws.Cells.Item(55, 1).Value = 'half'
$ws.Cells.Item(55, 2).Value = 'はん'
$ws.Cells.Item(56, 1).Value = 'half past two'
$ws.Cells.Item(56, 2).Value = 'にじはん'
$ws.Cells.Item(57, 1).Value = 'Japan'
$ws.Cells.Item(57, 2).Value = 'にほん'
$ws.Cells.Item(58, 1).Value = 'U.S.A.'
$ws.Cells.Item(58, 2).Value = 'アメリカ'
$ws.Cells.Item(59, 1).Value = '...language'
$ws.Cells.Item(59, 2).Value = '～ご'
$ws.Cells.Item(60, 1).Value = 'Japanese language'
$ws.Cells.Item(60, 2).Value = 'にほんご'
$ws.Cells.Item(61, 1).Value = '...years old'
$ws.Cells.Item(61, 2).Value = '～さい'
$ws.Cells.Item(62, 1).Value = 'telephone'
$ws.Cells.Item(62, 2).Value = 'でんわ'
$ws.Cells.Item(63, 1).Value = 'number...'
$ws.Cells.Item(63, 2).Value = '～ばん'
$ws.Cells.Item(64, 1).Value = 'number'
$ws.Cells.Item(64, 2).Value = 'ばんごう'
$ws.Cells.Item(65, 1).Value = 'name'
$ws.Cells.Item(65, 2).Value = 'なまえ'
$ws.Cells.Item(66, 1).Value = 'what'
$ws.Cells.Item(66, 2).Value = 'なん／なに'
$ws.Cells.Item(67, 1).Value = 'um...'
$ws.Cells.Item(67, 2).Value = 'あのう'
$ws.Cells.Item(68, 1).Value = 'yes'
$ws.Cells.Item(68, 2).Value = 'はい'
$ws.Cells.Item(69, 1).Value = 'That''s right.'
$ws.Cells.Item(69, 2).Value = 'そうです'
$ws.Cells.Item(70, 1).Value = 'I see.; Is that so?'
$ws.Cells.Item(70, 2).Value = 'そうですか'
$ws.Cells.Item(71, 1).Value = 'Britain'
$ws.Cells.Item(71, 2).Value = 'イギリス'
$ws.Cells.Item(72, 1).Value = 'Australia'
$ws.Cells.Item(72, 2).Value = 'オーストラリア'
$ws.Cells.Item(73, 1).Value = 'Korea'
$ws.Cells.Item(73, 2).Value = 'かんこく'
$ws.Cells.Item(74, 1).Value = 'Canada'
$ws.Cells.Item(74, 2).Value = 'カナダ'
$ws.Cells.Item(75, 1).Value = 'China'
$ws.Cells.Item(75, 2).Value = 'ちゅうごく'
$ws.Cells.Item(76, 1).Value = 'India'
$ws.Cells.Item(76, 2).Value = 'インド'
$ws.Cells.Item(77, 1).Value = 'Egypt'
$ws.Cells.Item(77, 2).Value = 'エジプト'
$ws.Cells.Item(78, 1).Value = 'Philippines'
$ws.Cells.Item(78, 2).Value = 'フィリピン'
$ws.Cells.Item(79, 1).Value = 'Japan'
$ws.Cells.Item(79, 2).Value = 'にほん'
$ws.Cells.Item(80, 1).Value = 'U.S.A.'
$ws.Cells.Item(80, 2).Value = 'アメリカ'
$ws.Cells.Item(81, 1).Value = 'Country'
$ws.Cells.Item(81, 2).Value = 'くに'
$ws.Cells.Item(82, 1).Value = 'Asian studies'
$ws.Cells.Item(82, 2).Value = 'アジアけんきゅう'
$ws.Cells.Item(83, 1).Value = 'economics'
$ws.Cells.Item(83, 2).Value = 'けいざい'
$ws.Cells.Item(84, 1).Value = 'engineering'
$ws.Cells.Item(84, 2).Value = 'こうがく'
$ws.Cells.Item(85, 1).Value = 'international relations'
$ws.Cells.Item(85, 2).Value = 'こくさいかんけい'
$ws.Cells.Item(86, 1).Value = 'computer'
$ws.Cells.Item(86, 2).Value = 'コンピューター'
$ws.Cells.Item(87, 1).Value = 'politics'
$ws.Cells.Item(87, 2).Value = 'せいじ'
$ws.Cells.Item(88, 1).Value = 'biology'
$ws.Cells.Item(88, 2).Value = 'せいぶつがく'
$ws.Cells.Item(89, 1).Value = 'business'
$ws.Cells.Item(89, 2).Value = 'ビジネス'
$ws.Cells.Item(90, 1).Value = 'literature'
$ws.Cells.Item(90, 2).Value = 'ぶんがく'
$ws.Cells.Item(91, 1).Value = 'history'
$ws.Cells.Item(91, 2).Value = 'れきし'
$ws.Cells.Item(92, 1).Value = 'major'
$ws.Cells.Item(92, 2).Value = 'せんこう'
$ws.Cells.Item(93, 1).Value = 'doctor'
$ws.Cells.Item(93, 2).Value = 'いしゃ'
$ws.Cells.Item(94, 1).Value = 'office worker'
$ws.Cells.Item(94, 2).Value = 'かいしゃいん'
$ws.Cells.Item(95, 1).Value = 'nurse'
$ws.Cells.Item(95, 2).Value = 'かんごし'
$ws.Cells.Item(96, 1).Value = 'high school student'
$ws.Cells.Item(96, 2).Value = 'こうこうせい'
$ws.Cells.Item(97, 1).Value = 'housewife'
$ws.Cells.Item(97, 2).Value = 'しゅふ'
$ws.Cells.Item(98, 1).Value = 'graduate student'
$ws.Cells.Item(98, 2).Value = 'だいがくいんせい'
$ws.Cells.Item(99, 1).Value = 'lawyer'
$ws.Cells.Item(99, 2).Value = 'べんごし'
$ws.Cells.Item(100, 1).Value = 'job; occupation'
$ws.Cells.Item(100, 2).Value = 'しごと'
$ws.Cells.Item(101, 1).Value = 'mother'
$ws.Cells.Item(101, 2).Value = 'おかあさん'
$ws.Cells.Item(102, 1).Value = 'father'
$ws.Cells.Item(102, 2).Value = 'おとうさん'
$ws.Cells.Item(103, 1).Value = 'older sister'
$ws.Cells.Item(103, 2).Value = 'おねえさん'
$ws.Cells.Item(104, 1).Value = 'older brother'
$ws.Cells.Item(104, 2).Value = 'おにいさん'
$ws.Cells.Item(105, 1).Value = 'younger sister'
$ws.Cells.Item(105, 2).Value = 'いもうと'
$ws.Cells.Item(106, 1).Value = 'younger brother'
$ws.Cells.Item(106, 2).Value = 'おとうと'
$ws.Cells.Item(107, 1).Value = 'family'
$ws.Cells.Item(107, 2).Value = 'かぞく'
$ws.Cells.Item(108, 1).Value = '...year student'
$ws.Cells.Item(108, 2).Value = '～ねんせい'
$ws.Cells.Item(109, 1).Value = 'first-year student'
$ws.Cells.Item(109, 2).Value = 'いちねんせい'
$ws.Cells.Item(110, 1).Value = 'second-year student'
$ws.Cells.Item(110, 2).Value = 'にねんせい'
$ws.Cells.Item(111, 1).Value = 'third-year student'
$ws.Cells.Item(111, 2).Value = 'さんねんせい'
$ws.Cells.Item(112, 1).Value = 'fourth-year student'
$ws.Cells.Item(112, 2).Value = 'よねんせい'
$ws.Cells.Item(113, 1).Value = 'fifth-year student'
$ws.Cells.Item(113, 2).Value = 'ごねんせい'
$ws.Cells.Item(114, 1).Value = 'sixth-year student'
$ws.Cells.Item(114, 2).Value = 'ろくねんせい'